# AFDP-3458: Add new MyDocuments module
#  - Add views privileges for new module
#  - Add creator as default assignee for PERSONAL repositories
#  - Deny read access to * participant
#
# This adds one new rule row (row 29) to the "Assignment Rules" table on
# Sheet1: "DocumentRepository - Default assignee" for DOC_REPO objects
# whose repositoryType is PERSONAL, assigning the object's creator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Best-effort: the source workbook also turned on iterative calculation
# (calcPr/@iterateDelta). Not every host surfaces this, so failures here
# are harmless and ignored.
try {
    $excel.Iteration = $true
    $excel.MaxChange = 0.0001
} catch { }

# The last existing data row is row 28 ("DocumentRepository - Default
# access"). Clone its formatting (borders, number formats, wrap, etc.)
# into the new row 29 so the new row looks consistent with the rest of
# the table.
$ws.Range("B28:H28").Copy()
$ws.Range("B29:H29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Column C on row 28 uses a plain (no-numberformat) style because its
# value ("DOC_REPO") happened to be entered with a leading quote on that
# row; row 29's C cell should use the boolean-styled format (matching
# column D), so pull that format from D28 instead.
$ws.Range("D28").Copy()
$ws.Range("C29").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row height used by the other wrapped-text rows in the table.
$ws.Rows("29").RowHeight = 30

# Fill in the new rule's cells. Order mirrors the order these new strings
# were first introduced in the authored workbook.
$ws.Range("B29").Value = "DocumentRepository – Default assignee"
$ws.Range("E29").Value = "repositoryType == 'PERSONAL'"
$ws.Range("H29").Value = "assignee, creator"
$ws.Range("D29").Value = "creator != null"
$ws.Range("C29").Value = "DOC_REPO"

# Move the on-screen selection/scroll position to the newly added row,
# matching where the author ended up after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D29").Select()
